# edit.ps1 - Apply "Ver 7 - Mid Level Contact" changes
# 1) Move the _GoBack bookmark from the end of the document (just before
#    "Sincerely,") to the very start of the body (right after the BODY
#    bookmark, replacing the gramStart/gramEnd proofing marks around "Dear ").
# 2) Split the "Error! Bookmark not defined." field-result run and the
#    "On behalf of everyone here at " run, adding spellStart/spellEnd
#    proofErr marks around "defined./On" (mirrors Word's automatic
#    spell-check marks that appear when the mid-level-contact name is
#    retyped).
# 3) Insert a new "SET SIGNATURE "DAVE"" field right after "Sincerely,"
#    (where the _GoBack bookmark used to start), before the relocated
#    _GoBack bookmark end.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$para1Xml = '<w:p w:rsidR="00946D1E" w:rsidRPr="00946D1E" w:rsidRDefault="00946D1E" w:rsidP="00D6148B"><w:pPr><w:spacing w:after="200"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="BODY"/><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">Dear </w:t></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="00ED17B0"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText>DOCVARIABLE var</w:instrText></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText xml:space="preserve">Greeting </w:instrText></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Error! Bookmark not defined.</w:t></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>,</w:t></w:r></w:p>'

$para2Xml = '<w:p w:rsidR="00A516EB" w:rsidRPr="00A516EB" w:rsidRDefault="00A516EB" w:rsidP="00D6148B"><w:pPr><w:spacing w:after="200"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r w:rsidRPr="00A516EB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">I am writing to express our sincere thanks for </w:t></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidR="00ED17B0"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText>DOCVARIABLE varReceipt</w:instrText></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">Error! Bookmark not </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>defined.</w:t></w:r><w:r w:rsidRPr="00946D1E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidRPr="00A516EB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>On</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> behalf of everyone here at </w:t></w:r><w:r w:rsidR="00640AF8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Partners In Health (</w:t></w:r><w:r w:rsidRPr="00A516EB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>PIH</w:t></w:r><w:r w:rsidR="00640AF8"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>)</w:t></w:r><w:r w:rsidRPr="00A516EB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>, I’d like to welcome you to our community—</w:t></w:r><w:r w:rsidR="00B829EE" w:rsidRPr="00B829EE"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00B829EE" w:rsidRPr="00B829EE"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>and tell you this donation helps us locate and treat people suffering from Ebola in West Africa.</w:t></w:r><w:r w:rsidRPr="00A516EB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$paraLastXml = '<w:p w:rsidR="00A516EB" w:rsidRPr="00A516EB" w:rsidRDefault="00A516EB" w:rsidP="00D6148B"><w:pPr><w:spacing w:after="200"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r w:rsidRPr="00A516EB"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:noProof/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Sincerely,</w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText>SET SIGNATURE "DAVE"</w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:instrText xml:space="preserve"> </w:instrText></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r><w:bookmarkEnd w:id="1"/></w:p>'

$paraCount = $d.Paragraphs.Count

$p1 = $d.Paragraphs(1).Range
$p1.InsertXML($para1Xml)

$p2 = $d.Paragraphs(2).Range
$p2.InsertXML($para2Xml)

$pLast = $d.Paragraphs($paraCount).Range
$pLast.InsertXML($paraLastXml)
